# Issue 864: Update orientation guide to reflect switch to TestNG and patch workflow.
#
# Two kinds of edits:
#  1) Every cached "datetimeFigureOut" field (slide master, all slide
#     layouts, and the notes master) gets its stale cached date text
#     refreshed from the old date to the new one (format depends on the
#     locale of the placeholder: en-US uses M/D/YYYY, en-SG uses D/M/YYYY).
#  2) On the one real slide, the code-block shape that reads
#     "     Java, JUnit" is updated to say "     Java, TestNG" instead,
#     reflecting the move from JUnit to TestNG.

$p = $ppt.ActivePresentation

# --- 1a. Slide master date placeholder -------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "12/17/2012") {
            $tr.Text = "5/18/2013"
        }
    }
}

# --- 1b. Every slide layout's date placeholder ------------------------
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "12/17/2012") {
                $tr.Text = "5/18/2013"
            }
        }
    }
}

# --- 1c. Notes master date placeholder (day/month/year order) ---------
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shp = $notesMaster.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "17/12/2012") {
            $tr.Text = "18/5/2013"
        }
    }
}

# --- 2. JUnit -> TestNG on the slide ----------------------------------
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $full = $shp.TextFrame.TextRange.Text
        $idx = $full.IndexOf("JUnit")
        if ($idx -ge 0) {
            $run = $shp.TextFrame.TextRange.Characters($idx + 1, 5)
            $run.Text = "TestNG"
        }
    }
}
